$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G2").Value = 0.1717705
$ws.Range("H2").Value = 0.343541
$ws.Range("I2").Value = 0.04063177891664595
$ws.Range("J2").Value = 0.02745976565347561
$ws.Range("M2").Value = 16.558025
$ws.Range("N2").Value = 33.11605
$ws.Range("O2").Value = 0.2047872600336892
$ws.Range("P2").Value = 0.155585835247525
$ws.Range("Q2").Value = 2.8441802332625
$ws.Range("R2").Value = 11.37672093305
$ws.Range("S2").Value = 0.008320870674634545
$ws.Range("T2").Value = 0.004272350574897302
$ws.Range("G3").Value = 0.1717705
$ws.Range("H3").Value = 0.343541
$ws.Range("I3").Value = 0.04063177891664595
$ws.Range("J3").Value = 0.02745976565347561
$ws.Range("M3").Value = 20.31779433333333
$ws.Range("O3").Value = 0.2512875437409564
$ws.Range("P3").Value = 0.2863712008291233
$ws.Range("Q3").Value = 3.489997691533833
$ws.Range("R3").Value = 20.939986149203
$ws.Range("S3").Value = 0.01021025992178954
$ws.Range("T3").Value = 0.007863686064672129
$ws.Range("G4").Value = 0.1717705
$ws.Range("H4").Value = 0.343541
$ws.Range("I4").Value = 0.04063177891664595
$ws.Range("J4").Value = 0.02745976565347561
$ws.Range("M4").Value = 10.405091
$ws.Range("N4").Value = 31.215273
$ws.Range("O4").Value = 0.1286886616182304
$ws.Range("P4").Value = 0.1466556042216543
$ws.Range("Q4").Value = 1.7872876836155
$ws.Range("R4").Value = 10.723726101693
$ws.Range("S4").Value = 0.005228849247950998
$ws.Range("T4").Value = 0.004027128523695495
$ws.Range("G5").Value = 0.1717705
$ws.Range("H5").Value = 0.343541
$ws.Range("I5").Value = 0.04063177891664595
$ws.Range("J5").Value = 0.02745976565347561
$ws.Range("M5").Value = 13.1587975
$ws.Range("N5").Value = 26.317595
$ws.Range("O5").Value = 0.1627461056112162
$ws.Range("P5").Value = 0.1236453320906656
$ws.Range("Q5").Value = 2.26029322597375
$ws.Range("R5").Value = 9.041172903894999
$ws.Range("S5").Value = 0.006612663782740052
$ws.Range("T5").Value = 0.003395271843355846
$ws.Range("G6").Value = 0.1717705
$ws.Range("H6").Value = 0.343541
$ws.Range("I6").Value = 0.04063177891664595
$ws.Range("J6").Value = 0.02745976565347561
$ws.Range("M6").Value = 7.303315666666667
$ws.Range("N6").Value = 21.909947
$ws.Range("O6").Value = 0.09032635260170116
$ws.Range("P6").Value = 0.1029373190408882
$ws.Range("Q6").Value = 1.254494183721167
$ws.Range("R6").Value = 7.526965102327001
$ws.Range("S6").Value = 0.003670120389259329
$ws.Range("T6").Value = 0.002826634657859842
$ws.Range("G7").Value = 0.1717705
$ws.Range("H7").Value = 0.343541
$ws.Range("I7").Value = 0.04063177891664595
$ws.Range("J7").Value = 0.02745976565347561
$ws.Range("M7").Value = 13.11173766666667
$ws.Range("N7").Value = 39.335213
$ws.Range("O7").Value = 0.1621640763942067
$ws.Range("P7").Value = 0.1848047085701436
$ws.Range("Q7").Value = 2.252209734872167
$ws.Range("R7").Value = 13.513258409233
$ws.Range("S7").Value = 0.006589014900271491
$ws.Range("T7").Value = 0.005074693988995
$ws.Range("I8").Value = 0.9593682210833541
$ws.Range("J8").Value = 0.9725402343465244
$ws.Range("M8").Value = 16.558025
$ws.Range("N8").Value = 33.11605
$ws.Range("O8").Value = 0.2047872600336892
$ws.Range("P8").Value = 0.155585835247525
$ws.Range("Q8").Value = 67.154729711025
$ws.Range("R8").Value = 402.92837826615
$ws.Range("S8").Value = 0.1964663893590546
$ws.Range("T8").Value = 0.1513134846726277
$ws.Range("I9").Value = 0.9593682210833541
$ws.Range("J9").Value = 0.9725402343465244
$ws.Range("M9").Value = 20.31779433333333
$ws.Range("O9").Value = 0.2512875437409564
$ws.Range("P9").Value = 0.2863712008291233
$ws.Range("Q9").Value = 82.40330515138099
$ws.Range("R9").Value = 741.6297463624289
$ws.Range("S9").Value = 0.2410772838191669
$ws.Range("T9").Value = 0.2785075147644512
$ws.Range("I10").Value = 0.9593682210833541
$ws.Range("J10").Value = 0.9725402343465244
$ws.Range("M10").Value = 10.405091
$ws.Range("N10").Value = 31.215273
$ws.Range("O10").Value = 0.1286886616182304
$ws.Range("P10").Value = 0.1466556042216543
$ws.Range("Q10").Value = 42.20014607561099
$ws.Range("R10").Value = 379.801314680499
$ws.Range("S10").Value = 0.1234598123702794
$ws.Range("T10").Value = 0.1426284756979588
$ws.Range("I11").Value = 0.9593682210833541
$ws.Range("J11").Value = 0.9725402343465244
$ws.Range("M11").Value = 13.1587975
$ws.Range("N11").Value = 26.317595
$ws.Range("O11").Value = 0.1627461056112162
$ws.Range("P11").Value = 0.1236453320906656
$ws.Range("Q11").Value = 53.36841135549749
$ws.Range("R11").Value = 320.210468132985
$ws.Range("S11").Value = 0.1561334418284762
$ws.Range("T11").Value = 0.1202500602473098
$ws.Range("I12").Value = 0.9593682210833541
$ws.Range("J12").Value = 0.9725402343465244
$ws.Range("M12").Value = 7.303315666666667
$ws.Range("N12").Value = 21.909947
$ws.Range("O12").Value = 0.09032635260170116
$ws.Range("P12").Value = 0.1029373190408882
$ws.Range("Q12").Value = 29.620210718929
$ws.Range("R12").Value = 266.581896470361
$ws.Range("S12").Value = 0.08665623221244183
$ws.Range("T12").Value = 0.1001106843830283
$ws.Range("I13").Value = 0.9593682210833541
$ws.Range("J13").Value = 0.9725402343465244
$ws.Range("M13").Value = 13.11173766666667
$ws.Range("N13").Value = 39.335213
$ws.Range("O13").Value = 0.1621640763942067
$ws.Range("P13").Value = 0.1848047085701436
$ws.Range("Q13").Value = 53.17754980119101
$ws.Range("R13").Value = 478.597948210719
$ws.Range("S13").Value = 0.1555750614939352
$ws.Range("T13").Value = 0.1797300145811486
